$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-10-28 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-29 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("68-49=19", $true, $false, $false, $false, $false, $true, 1, $false, "3+78=81", 2) | Out-Null
$d.Content.Find.Execute("78+21=99", $true, $false, $false, $false, $false, $true, 1, $false, "94-42=52", 2) | Out-Null
$d.Content.Find.Execute("57+30=87", $true, $false, $false, $false, $false, $true, 1, $false, "23+38=61", 2) | Out-Null
$d.Content.Find.Execute("25+0=25", $true, $false, $false, $false, $false, $true, 1, $false, "53-40=13", 2) | Out-Null
$d.Content.Find.Execute("76-39=37", $true, $false, $false, $false, $false, $true, 1, $false, "15+11=26", 2) | Out-Null
$d.Content.Find.Execute("61-10=51", $true, $false, $false, $false, $false, $true, 1, $false, "32+47=79", 2) | Out-Null
$d.Content.Find.Execute("19+20=39", $true, $false, $false, $false, $false, $true, 1, $false, "81-46=35", 2) | Out-Null
$d.Content.Find.Execute("24+28=52", $true, $false, $false, $false, $false, $true, 1, $false, "57-19=38", 2) | Out-Null
$d.Content.Find.Execute("64-35=29", $true, $false, $false, $false, $false, $true, 1, $false, "53-11=42", 2) | Out-Null
$d.Content.Find.Execute("0+74=74", $true, $false, $false, $false, $false, $true, 1, $false, "35+24=59", 2) | Out-Null
$d.Content.Find.Execute("19+25=44", $true, $false, $false, $false, $false, $true, 1, $false, "78+11=89", 2) | Out-Null
$d.Content.Find.Execute("55-11=44", $true, $false, $false, $false, $false, $true, 1, $false, "18+23=41", 2) | Out-Null
$d.Content.Find.Execute("45+39=84", $true, $false, $false, $false, $false, $true, 1, $false, "83-1=82", 2) | Out-Null
$d.Content.Find.Execute("2+19=21", $true, $false, $false, $false, $false, $true, 1, $false, "61+36=97", 2) | Out-Null
$d.Content.Find.Execute("11+45=56", $true, $false, $false, $false, $false, $true, 1, $false, "88+10=98", 2) | Out-Null
$d.Content.Find.Execute("81-33=48", $true, $false, $false, $false, $false, $true, 1, $false, "58+19=77", 2) | Out-Null
$d.Content.Find.Execute("35-5=30", $true, $false, $false, $false, $false, $true, 1, $false, "34+43=77", 2) | Out-Null
$d.Content.Find.Execute("41+8=49", $true, $false, $false, $false, $false, $true, 1, $false, "91-28=63", 2) | Out-Null
$d.Content.Find.Execute("56-52=4", $true, $false, $false, $false, $false, $true, 1, $false, "2+52=54", 2) | Out-Null
$d.Content.Find.Execute("75-9=66", $true, $false, $false, $false, $false, $true, 1, $false, "84-1=83", 2) | Out-Null
$d.Content.Find.Execute("51-29=22", $true, $false, $false, $false, $false, $true, 1, $false, "17+81=98", 2) | Out-Null
$d.Content.Find.Execute("32+6=38", $true, $false, $false, $false, $false, $true, 1, $false, "47-34=13", 2) | Out-Null
$d.Content.Find.Execute("89-76=13", $true, $false, $false, $false, $false, $true, 1, $false, "47-23=24", 2) | Out-Null
$d.Content.Find.Execute("88-61=27", $true, $false, $false, $false, $false, $true, 1, $false, "21-20=1", 2) | Out-Null
$d.Content.Find.Execute("48+33=81", $true, $false, $false, $false, $false, $true, 1, $false, "85-66=19", 2) | Out-Null
$d.Content.Find.Execute("84-32=52", $true, $false, $false, $false, $false, $true, 1, $false, "55-9=46", 2) | Out-Null
$d.Content.Find.Execute("70-46=24", $true, $false, $false, $false, $false, $true, 1, $false, "17+37=54", 2) | Out-Null
$d.Content.Find.Execute("96-59=37", $true, $false, $false, $false, $false, $true, 1, $false, "51-23=28", 2) | Out-Null
$d.Content.Find.Execute("66+29=95", $true, $false, $false, $false, $false, $true, 1, $false, "88+8=96", 2) | Out-Null
$d.Content.Find.Execute("22+38=60", $true, $false, $false, $false, $false, $true, 1, $false, "98-32=66", 2) | Out-Null
$d.Content.Find.Execute("87-41=46", $true, $false, $false, $false, $false, $true, 1, $false, "57+18=75", 2) | Out-Null
$d.Content.Find.Execute("75+8=83", $true, $false, $false, $false, $false, $true, 1, $false, "86+0=86", 2) | Out-Null
$d.Content.Find.Execute("98-78=20", $true, $false, $false, $false, $false, $true, 1, $false, "5-3=2", 2) | Out-Null
$d.Content.Find.Execute("82-35=47", $true, $false, $false, $false, $false, $true, 1, $false, "43-9=34", 2) | Out-Null
$d.Content.Find.Execute("3+81=84", $true, $false, $false, $false, $false, $true, 1, $false, "64+16=80", 2) | Out-Null
$d.Content.Find.Execute("90-82=8", $true, $false, $false, $false, $false, $true, 1, $false, "91-22=69", 2) | Out-Null
$d.Content.Find.Execute("82-81=1", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=81", 2) | Out-Null
$d.Content.Find.Execute("85-33=52", $true, $false, $false, $false, $false, $true, 1, $false, "70-38=32", 2) | Out-Null
$d.Content.Find.Execute("88-50=38", $true, $false, $false, $false, $false, $true, 1, $false, "81-51=30", 2) | Out-Null
$d.Content.Find.Execute("53+33=86", $true, $false, $false, $false, $false, $true, 1, $false, "71+8=79", 2) | Out-Null
$d.Content.Find.Execute("33+22=55", $true, $false, $false, $false, $false, $true, 1, $false, "16+11=27", 2) | Out-Null
$d.Content.Find.Execute("23+55=78", $true, $false, $false, $false, $false, $true, 1, $false, "27+71=98", 2) | Out-Null
$d.Content.Find.Execute("23+28=51", $true, $false, $false, $false, $false, $true, 1, $false, "12+81=93", 2) | Out-Null
$d.Content.Find.Execute("77+5=82", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=83", 2) | Out-Null
$d.Content.Find.Execute("51-39=12", $true, $false, $false, $false, $false, $true, 1, $false, "76-43=33", 2) | Out-Null
$d.Content.Find.Execute("51+24=75", $true, $false, $false, $false, $false, $true, 1, $false, "5+56=61", 2) | Out-Null
$d.Content.Find.Execute("17-9=8", $true, $false, $false, $false, $false, $true, 1, $false, "74-20=54", 2) | Out-Null
$d.Content.Find.Execute("43-30=13", $true, $false, $false, $false, $false, $true, 1, $false, "4+69=73", 2) | Out-Null
$d.Content.Find.Execute("35+11=46", $true, $false, $false, $false, $false, $true, 1, $false, "39-15=24", 2) | Out-Null
$d.Content.Find.Execute("66-56=10", $true, $false, $false, $false, $false, $true, 1, $false, "83-50=33", 2) | Out-Null
$d.Content.Find.Execute("0+1=1", $true, $false, $false, $false, $false, $true, 1, $false, "33+30=63", 2) | Out-Null
$d.Content.Find.Execute("15+80=95", $true, $false, $false, $false, $false, $true, 1, $false, "33+36=69", 2) | Out-Null
$d.Content.Find.Execute("34-3=31", $true, $false, $false, $false, $false, $true, 1, $false, "45-16=29", 2) | Out-Null
$d.Content.Find.Execute("58-54=4", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=37", 2) | Out-Null
$d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "55-34=21", 2) | Out-Null
$d.Content.Find.Execute("17+18=35", $true, $false, $false, $false, $false, $true, 1, $false, "9+41=50", 2) | Out-Null
$d.Content.Find.Execute("32+52=84", $true, $false, $false, $false, $false, $true, 1, $false, "30+53=83", 2) | Out-Null
$d.Content.Find.Execute("68-67=1", $true, $false, $false, $false, $false, $true, 1, $false, "48+44=92", 2) | Out-Null
$d.Content.Find.Execute("62-15=47", $true, $false, $false, $false, $false, $true, 1, $false, "91-3=88", 2) | Out-Null
$d.Content.Find.Execute("82-29=53", $true, $false, $false, $false, $false, $true, 1, $false, "40-4=36", 2) | Out-Null
$d.Content.Find.Execute("35-15=20", $true, $false, $false, $false, $false, $true, 1, $false, "8+33=41", 2) | Out-Null
$d.Content.Find.Execute("90-31=59", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=8", 2) | Out-Null
$d.Content.Find.Execute("10+89=99", $true, $false, $false, $false, $false, $true, 1, $false, "65+28=93", 2) | Out-Null
$d.Content.Find.Execute("6+63=69", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=25", 2) | Out-Null
$d.Content.Find.Execute("36+31=67", $true, $false, $false, $false, $false, $true, 1, $false, "79-52=27", 2) | Out-Null
$d.Content.Find.Execute("1+18=19", $true, $false, $false, $false, $false, $true, 1, $false, "41-28=13", 2) | Out-Null
$d.Content.Find.Execute("78+10=88", $true, $false, $false, $false, $false, $true, 1, $false, "99-24=75", 2) | Out-Null
$d.Content.Find.Execute("3+45=48", $true, $false, $false, $false, $false, $true, 1, $false, "11+42=53", 2) | Out-Null
$d.Content.Find.Execute("34+65=99", $true, $false, $false, $false, $false, $true, 1, $false, "1+39=40", 2) | Out-Null
$d.Content.Find.Execute("22+42=64", $true, $false, $false, $false, $false, $true, 1, $false, "43+50=93", 2) | Out-Null
$d.Content.Find.Execute("7+14=21", $true, $false, $false, $false, $false, $true, 1, $false, "80+11=91", 2) | Out-Null
$d.Content.Find.Execute("53-5=48", $true, $false, $false, $false, $false, $true, 1, $false, "98-80=18", 2) | Out-Null
$d.Content.Find.Execute("3+51=54", $true, $false, $false, $false, $false, $true, 1, $false, "87-64=23", 2) | Out-Null
$d.Content.Find.Execute("31+20=51", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=66", 2) | Out-Null
$d.Content.Find.Execute("58+38=96", $true, $false, $false, $false, $false, $true, 1, $false, "4+65=69", 2) | Out-Null
$d.Content.Find.Execute("32+41=73", $true, $false, $false, $false, $false, $true, 1, $false, "28+16=44", 2) | Out-Null
$d.Content.Find.Execute("36-36=0", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=40", 2) | Out-Null
$d.Content.Find.Execute("14+8=22", $true, $false, $false, $false, $false, $true, 1, $false, "14-4=10", 2) | Out-Null
$d.Content.Find.Execute("35-33=2", $true, $false, $false, $false, $false, $true, 1, $false, "40-10=30", 2) | Out-Null
$d.Content.Find.Execute("7+57=64", $true, $false, $false, $false, $false, $true, 1, $false, "21+40=61", 2) | Out-Null
$d.Content.Find.Execute("39-20=19", $true, $false, $false, $false, $false, $true, 1, $false, "67-59=8", 2) | Out-Null
$d.Content.Find.Execute("12+56=68", $true, $false, $false, $false, $false, $true, 1, $false, "36+9=45", 2) | Out-Null
$d.Content.Find.Execute("87-60=27", $true, $false, $false, $false, $false, $true, 1, $false, "46+14=60", 2) | Out-Null
$d.Content.Find.Execute("34-2=32", $true, $false, $false, $false, $false, $true, 1, $false, "96-38=58", 2) | Out-Null
$d.Content.Find.Execute("87-12=75", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=29", 2) | Out-Null
$d.Content.Find.Execute("29+47=76", $true, $false, $false, $false, $false, $true, 1, $false, "19+80=99", 2) | Out-Null
$d.Content.Find.Execute("21+45=66", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=44", 2) | Out-Null
$d.Content.Find.Execute("25+18=43", $true, $false, $false, $false, $false, $true, 1, $false, "30+51=81", 2) | Out-Null
$d.Content.Find.Execute("58+11=69", $true, $false, $false, $false, $false, $true, 1, $false, "36+8=44", 2) | Out-Null
$d.Content.Find.Execute("44+8=52", $true, $false, $false, $false, $false, $true, 1, $false, "78-75=3", 2) | Out-Null
$d.Content.Find.Execute("34+0=34", $true, $false, $false, $false, $false, $true, 1, $false, "97-68=29", 2) | Out-Null
$d.Content.Find.Execute("51-13=38", $true, $false, $false, $false, $false, $true, 1, $false, "87-75=12", 2) | Out-Null
$d.Content.Find.Execute("23+63=86", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 2) | Out-Null
$d.Content.Find.Execute("69+19=88", $true, $false, $false, $false, $false, $true, 1, $false, "75-59=16", 2) | Out-Null
$d.Content.Find.Execute("66-6=60", $true, $false, $false, $false, $false, $true, 1, $false, "61+2=63", 2) | Out-Null
$d.Content.Find.Execute("37-1=36", $true, $false, $false, $false, $false, $true, 1, $false, "16+80=96", 2) | Out-Null
$d.Content.Find.Execute("90-66=24", $true, $false, $false, $false, $false, $true, 1, $false, "37+52=89", 2) | Out-Null
$d.Content.Find.Execute("96+3=99", $true, $false, $false, $false, $false, $true, 1, $false, "66+20=86", 2) | Out-Null
$d.Content.Find.Execute("50+1=51", $true, $false, $false, $false, $false, $true, 1, $false, "74-8=66", 2) | Out-Null
$d.Content.Find.Execute("31+24=55", $true, $false, $false, $false, $false, $true, 1, $false, "98-62=36", 2) | Out-Null
